$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 310, pushing the existing rows 310:420 down to 312:422
$ws.Rows("310:311").Insert()

# New row 310 (Primera)
$ws.Range("A310").Value = 11
$ws.Range("B310").Value = "Vega Monumental Concepción"
$ws.Range("C310").Value = "Bíobío"
$ws.Range("D310").Value = 44988
$ws.Range("E310").Value = 8
$ws.Range("F310").Value = 100112017
$ws.Range("G310").Value = "Apio"
$ws.Range("H310").Value = "Americana (o)"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 220
$ws.Range("K310").Value = 8500
$ws.Range("L310").Value = 9000
$ws.Range("M310").Value = 8773
$ws.Range("N310").Value = '$/docena de matas'
$ws.Range("O310").Value = "Región de Coquimbo"
$ws.Range("P310").Value = 1462
$ws.Range("Q310").Value = 6
$ws.Range("R310").Value = "Hortaliza"

# New row 311 (Segunda)
$ws.Range("A311").Value = 11
$ws.Range("B311").Value = "Vega Monumental Concepción"
$ws.Range("C311").Value = "Bíobío"
$ws.Range("D311").Value = 44988
$ws.Range("E311").Value = 8
$ws.Range("F311").Value = 100112017
$ws.Range("G311").Value = "Apio"
$ws.Range("H311").Value = "Americana (o)"
$ws.Range("I311").Value = "Segunda"
$ws.Range("J311").Value = 150
$ws.Range("K311").Value = 6500
$ws.Range("L311").Value = 6500
$ws.Range("M311").Value = 6500
$ws.Range("N311").Value = '$/docena de matas'
$ws.Range("O311").Value = "Región de Coquimbo"
$ws.Range("P311").Value = 1083
$ws.Range("Q311").Value = 6
$ws.Range("R311").Value = "Hortaliza"
